$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.001655459403992
$ws.Range("B1").Value = 1.935990691184998
$ws.Range("C1").Value = 3.068731069564819
$ws.Range("D1").Value = 3.752788305282593
$ws.Range("E1").Value = 1.637585282325745
